# posts.xlsx: remove the post row for "「ようこそ」اهلا وسهلا ..." (previously
# row 765). Excel's native row delete shifts every row beneath it up by one,
# so what was row 766 ("「まぁ何というサプライズでしょう」...") becomes the
# new row 765, and so on down through the former row 774, which becomes the
# new last row 773 — matching the updated dimension A1:C773.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(765).Delete()
